$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for author "Mart\'inez" (A. Martinez) - row 67.
# Deleting the row shifts every row below it up by one.
$ws.Rows.Item(67).Delete() | Out-Null

# Remove the 2nd affiliation (Institution2/Address2, columns F:G) for
# "Palmeiro" (B. Palmeiro), which after the row deletion above now sits
# on row 82 (was row 83).
$ws.Range("F82:G82").ClearContents() | Out-Null

# Match the author's final selection/scroll position.
$ws.Range("A68").Select() | Out-Null
